$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G (K)
$gValues = @{
    2 = 1
    3 = 5
    4 = 6
    5 = 6
    6 = 3
    7 = 3
    8 = 6
    9 = 5
    10 = 6
    11 = 9
    12 = 6
    13 = 4
    14 = 5
    15 = 6
    16 = 3
    17 = 7
    18 = 4
    19 = 4
    20 = 4
    21 = 4
    22 = 6
    23 = 2
    24 = 3
    25 = 9
    26 = 6
    27 = 4
    28 = 8
    29 = 6
    30 = 5
    31 = 5
    32 = 4
    33 = 4
    34 = 9
    35 = 3
    36 = 4
    37 = 1
    38 = 3
    39 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

